$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'96.818.85"
$ws.Range('E2').Value = '  -0.95%  '

$ws.Range('D3').Value = "'3.339.17"
$ws.Range('E3').Value = '  -1.97%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = "'250.70"
$ws.Range('E5').Value = '  -1.83%  '

$ws.Range('D6').Value = "'657.51"
$ws.Range('E6').Value = '  +1.09%  '

$ws.Range('E7').Value = '  -3.30%  '

$ws.Range('D8').Value = "'0.425"
$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('D9').Value = "'0.999"
$ws.Range('E9').Value = '  +0.06%  '

$ws.Range('E10').Value = '  -4.53%  '

$ws.Range('D11').Value = "'3.336.06"
$ws.Range('E11').Value = '  -1.96%  '

$ws.Range('E12').Value = '  -2.21%  '

$ws.Range('D13').Value = "'40.57"
$ws.Range('E13').Value = '  -1.66%  '

$ws.Range('D14').Value = "'96.589.03"
$ws.Range('E14').Value = '  -0.84%  '

$ws.Range('E15').Value = '  -2.39%  '

$ws.Range('E16').Value = '  -2.01%  '

$ws.Range('D17').Value = "'3.967.56"
$ws.Range('E17').Value = '  -1.61%  '

$ws.Range('D18').Value = "'8.72"
$ws.Range('E18').Value = '  +2.70%  '

$ws.Range('D19').Value = "'3.343.55"
$ws.Range('E19').Value = '  -1.50%  '

$ws.Range('D20').Value = "'0.564"
$ws.Range('E20').Value = '  +15.42%  '

$ws.Range('D21').Value = "'17.38"
$ws.Range('E21').Value = '  +0.58%  '

$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = "'507.77"
$ws.Range('E22').Value = '  +1.70%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = "'10.60"
$ws.Range('E23').Value = '  -0.91%  '

$ws.Range('E24').Value = '  -1.82%  '

$ws.Range('D25').Value = "'0.0000199"
$ws.Range('E25').Value = '  -2.55%  '

$ws.Range('D26').Value = "'6.56"
$ws.Range('E26').Value = '  +6.98%  '

$ws.Range('D27').Value = "'96.52"
$ws.Range('E27').Value = '  -1.77%  '

$ws.Range('D28').Value = "'12.13"
$ws.Range('E28').Value = '  -3.45%  '

$ws.Range('E29').Value = '  -4.41%  '

$ws.Range('E30').Value = '  +0.34%  '

$ws.Range('D31').Value = "'11.18"
$ws.Range('E31').Value = '  -0.32%  '

$ws.Range('E32').Value = '  -6.13%  '

$ws.Range('E33').Value = '  +12.78%  '

$ws.Range('E34').Value = '  +0.33%  '

$ws.Range('D35').Value = "'0.552"
$ws.Range('E35').Value = '  -2.56%  '

$ws.Range('D36').Value = "'28.33"
$ws.Range('E36').Value = '  -4.01%  '

$ws.Range('D37').Value = "'7.86"
$ws.Range('E37').Value = '  +2.40%  '

$ws.Range('D38').Value = "'1.49"
$ws.Range('E38').Value = '  +6.32%  '

$ws.Range('E39').Value = '  -0.38%  '

$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('D41').Value = "'508.71"
$ws.Range('E41').Value = '  -0.52%  '

$ws.Range('E42').Value = '  -1.40%  '

$ws.Range('D43').Value = "'0.0434"
$ws.Range('E43').Value = '  +4.89%  '

$ws.Range('E44').Value = '  -2.67%  '

$ws.Range('D45').Value = "'3.64"
$ws.Range('E45').Value = '  -1.09%  '

$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').Value = "'1.67"
$ws.Range('E46').Value = '  +7.32%  '

$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = "'5.56"
$ws.Range('E47').Value = '  +1.96%  '

$ws.Range('D48').Value = "'8.46"
$ws.Range('E48').Value = '  +3.62%  '

$ws.Range('D49').Value = "'53.24"
$ws.Range('E49').Value = '  +3.56%  '

$ws.Range('E50').Value = '  -2.87%  '

$ws.Range('E51').Value = '  +0.85%  '
